$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '53.931.01'
$ws.Range("E2").Value = '  -3.93%  '
$ws.Range("D3").Value = '2.257.77'
$ws.Range("E3").Value = '  -4.55%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '488.27'
$ws.Range("E5").Value = '  -1.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '125.98'
$ws.Range("E6").Value = '  -2.91%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.519'
$ws.Range("E8").Value = '  -4.16%  '
$ws.Range("D9").Value = '2.260.99'
$ws.Range("E9").Value = '  -4.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0919'
$ws.Range("E10").Value = '  -5.05%  '
$ws.Range("E11").Value = '  -1.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.75'
$ws.Range("E12").Value = '  +3.48%  '
$ws.Range("E13").Value = '  -2.89%  '
$ws.Range("D14").Value = '2.661.06'
$ws.Range("E14").Value = '  -4.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.96'
$ws.Range("E15").Value = '  -2.04%  '
$ws.Range("D16").Value = '53.883.72'
$ws.Range("E16").Value = '  -3.96%  '
$ws.Range("E17").Value = '  -2.58%  '
$ws.Range("D18").Value = '2.245.87'
$ws.Range("E18").Value = '  -6.13%  '
$ws.Range("E19").Value = '  -3.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.96'
$ws.Range("E20").Value = '  -0.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '300.61'
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("E22").Value = '  -1.96%  '
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.83'
$ws.Range("E24").Value = '  -1.80%  '
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("E27").Value = '  -2.66%  '
$ws.Range("E28").Value = '  -2.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '169.80'
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("D30").Value = '0.0₃0688'
$ws.Range("E30").Value = '  -2.79%  '
$ws.Range("E31").Value = '  -1.98%  '
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.72'
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("E35").Value = '  -1.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.39'
$ws.Range("E36").Value = '  -0.72%  '
$ws.Range("E37").Value = '  -0.91%  '
$ws.Range("E38").Value = '  +5.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.58'
$ws.Range("E39").Value = '  -4.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35.74'
$ws.Range("E40").Value = '  -0.64%  '
$ws.Range("E41").Value = '  -0.80%  '
$ws.Range("E42").Value = '  -2.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.28'
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '122.12'
$ws.Range("E44").Value = '  -6.26%  '
$ws.Range("E45").Value = '  -1.47%  '
$ws.Range("E46").Value = '  -2.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.537'
$ws.Range("E47").Value = '  -4.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '237.48'
$ws.Range("E48").Value = '  -0.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0471'
$ws.Range("E49").Value = '  -1.86%  '
$ws.Range("E50").Value = '  -2.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.19'
$ws.Range("E51").Value = '  -3.19%  '
